$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update changed cell values for rows 2-27 based on the new data
# Row 2
$ws.Cells.Item(2, 3).Value = 3
$ws.Cells.Item(2, 6).Value = 109
$ws.Cells.Item(2, 8).Value = "living_rooms"
$ws.Cells.Item(2, 12).Value = "stimuli/img_95hiv.png"
$ws.Cells.Item(2, 13).Value = 84.04545454545455
$ws.Cells.Item(2, 14).Value = 67.31818181818181
$ws.Cells.Item(2, 15).Value = 75.68181818181819
$ws.Cells.Item(2, 16).Value = 44
$ws.Cells.Item(2, 17).Value = 9
$ws.Cells.Item(2, 18).Value = 9
$ws.Cells.Item(2, 19).Value = 9
# Row 3
$ws.Cells.Item(3, 3).Value = 3
$ws.Cells.Item(3, 6).Value = 110
$ws.Cells.Item(3, 8).Value = "living_rooms"
$ws.Cells.Item(3, 12).Value = "stimuli/img_6a0hu.png"
$ws.Cells.Item(3, 13).Value = 61.275
$ws.Cells.Item(3, 14).Value = 42.025
$ws.Cells.Item(3, 15).Value = 51.65
$ws.Cells.Item(3, 16).Value = 40
$ws.Cells.Item(3, 17).Value = 4
$ws.Cells.Item(3, 18).Value = 4
$ws.Cells.Item(3, 19).Value = 4
# Row 4
$ws.Cells.Item(4, 3).Value = 3
$ws.Cells.Item(4, 6).Value = 111
$ws.Cells.Item(4, 8).Value = "kitchens"
$ws.Cells.Item(4, 12).Value = "stimuli/img_gztbt.png"
$ws.Cells.Item(4, 13).Value = 55.06451612903226
$ws.Cells.Item(4, 14).Value = 26.09677419354839
$ws.Cells.Item(4, 15).Value = 40.58064516129032
$ws.Cells.Item(4, 16).Value = 31
$ws.Cells.Item(4, 17).Value = 2
$ws.Cells.Item(4, 18).Value = 2
$ws.Cells.Item(4, 19).Value = 2
# Row 5
$ws.Cells.Item(5, 3).Value = 3
$ws.Cells.Item(5, 6).Value = 112
$ws.Cells.Item(5, 9).Value = "target"
$ws.Cells.Item(5, 11).Value = "j"
$ws.Cells.Item(5, 12).Value = "stimuli/img_tujn3.png"
$ws.Cells.Item(5, 13).Value = 81.4090909090909
$ws.Cells.Item(5, 14).Value = 62.52272727272727
$ws.Cells.Item(5, 15).Value = 71.9659090909091
$ws.Cells.Item(5, 16).Value = 44
$ws.Cells.Item(5, 17).Value = 8
$ws.Cells.Item(5, 18).Value = 8
$ws.Cells.Item(5, 19).Value = 8
# Row 6
$ws.Cells.Item(6, 3).Value = 3
$ws.Cells.Item(6, 6).Value = 113
$ws.Cells.Item(6, 9).Value = "distractor"
$ws.Cells.Item(6, 11).Value = "f"
$ws.Cells.Item(6, 12).Value = "stimuli/img_pt3d7.png"
$ws.Cells.Item(6, 13).Value = 65.08571428571429
$ws.Cells.Item(6, 14).Value = 44.65714285714286
$ws.Cells.Item(6, 15).Value = 54.87142857142857
$ws.Cells.Item(6, 16).Value = 35
$ws.Cells.Item(6, 17).Value = 4
$ws.Cells.Item(6, 18).Value = 4
$ws.Cells.Item(6, 19).Value = 4
# Row 7
$ws.Cells.Item(7, 3).Value = 3
$ws.Cells.Item(7, 6).Value = 114
$ws.Cells.Item(7, 8).Value = "living_rooms"
$ws.Cells.Item(7, 12).Value = "stimuli/img_cehin.png"
$ws.Cells.Item(7, 13).Value = 78.86363636363636
$ws.Cells.Item(7, 14).Value = 60.02272727272727
$ws.Cells.Item(7, 15).Value = 69.44318181818181
$ws.Cells.Item(7, 16).Value = 44
$ws.Cells.Item(7, 17).Value = 7
$ws.Cells.Item(7, 18).Value = 7
$ws.Cells.Item(7, 19).Value = 7
# Row 8
$ws.Cells.Item(8, 3).Value = 3
$ws.Cells.Item(8, 6).Value = 115
$ws.Cells.Item(8, 8).Value = "living_rooms"
$ws.Cells.Item(8, 12).Value = "stimuli/img_kost0.png"
$ws.Cells.Item(8, 13).Value = 63.09090909090909
$ws.Cells.Item(8, 14).Value = 42.77272727272727
$ws.Cells.Item(8, 15).Value = 52.93181818181819
$ws.Cells.Item(8, 16).Value = 44
$ws.Cells.Item(8, 17).Value = 5
$ws.Cells.Item(8, 18).Value = 5
$ws.Cells.Item(8, 19).Value = 5
# Row 9
$ws.Cells.Item(9, 3).Value = 3
$ws.Cells.Item(9, 6).Value = 116
$ws.Cells.Item(9, 8).Value = "bedrooms"
$ws.Cells.Item(9, 9).Value = "distractor"
$ws.Cells.Item(9, 11).Value = "f"
$ws.Cells.Item(9, 12).Value = "stimuli/img_die1d.png"
$ws.Cells.Item(9, 13).Value = 75.42857142857143
$ws.Cells.Item(9, 14).Value = 53.30952380952381
$ws.Cells.Item(9, 15).Value = 64.36904761904762
$ws.Cells.Item(9, 16).Value = 42
$ws.Cells.Item(9, 17).Value = 6
$ws.Cells.Item(9, 18).Value = 6
$ws.Cells.Item(9, 19).Value = 6
# Row 10
$ws.Cells.Item(10, 3).Value = 3
$ws.Cells.Item(10, 6).Value = 117
$ws.Cells.Item(10, 8).Value = "living_rooms"
$ws.Cells.Item(10, 12).Value = "stimuli/img_0kqc0.png"
$ws.Cells.Item(10, 13).Value = 43.74468085106383
$ws.Cells.Item(10, 14).Value = 27.14893617021277
$ws.Cells.Item(10, 15).Value = 35.4468085106383
$ws.Cells.Item(10, 16).Value = 47
$ws.Cells.Item(10, 17).Value = 2
$ws.Cells.Item(10, 18).Value = 2
$ws.Cells.Item(10, 19).Value = 2
# Row 11
$ws.Cells.Item(11, 3).Value = 3
$ws.Cells.Item(11, 6).Value = 118
$ws.Cells.Item(11, 8).Value = "kitchens"
$ws.Cells.Item(11, 12).Value = "stimuli/img_7wquy.png"
$ws.Cells.Item(11, 13).Value = 50.59375
$ws.Cells.Item(11, 14).Value = 30.59375
$ws.Cells.Item(11, 15).Value = 40.59375
$ws.Cells.Item(11, 16).Value = 32
$ws.Cells.Item(11, 17).Value = 2
$ws.Cells.Item(11, 18).Value = 2
$ws.Cells.Item(11, 19).Value = 2
# Row 12
$ws.Cells.Item(12, 3).Value = 3
$ws.Cells.Item(12, 6).Value = 119
$ws.Cells.Item(12, 8).Value = "living_rooms"
$ws.Cells.Item(12, 12).Value = "stimuli/img_wz6x5.png"
$ws.Cells.Item(12, 13).Value = 68.3695652173913
$ws.Cells.Item(12, 14).Value = 48.47826086956522
$ws.Cells.Item(12, 15).Value = 58.42391304347826
$ws.Cells.Item(12, 16).Value = 46
$ws.Cells.Item(12, 17).Value = 5
$ws.Cells.Item(12, 18).Value = 5
$ws.Cells.Item(12, 19).Value = 5
# Row 13
$ws.Cells.Item(13, 3).Value = 3
$ws.Cells.Item(13, 6).Value = 120
$ws.Cells.Item(13, 8).Value = "living_rooms"
$ws.Cells.Item(13, 12).Value = "stimuli/img_xu1p3.png"
$ws.Cells.Item(13, 13).Value = 75.27659574468085
$ws.Cells.Item(13, 14).Value = 56.68085106382978
$ws.Cells.Item(13, 15).Value = 65.97872340425532
$ws.Cells.Item(13, 16).Value = 47
$ws.Cells.Item(13, 17).Value = 7
$ws.Cells.Item(13, 18).Value = 7
$ws.Cells.Item(13, 19).Value = 7
# Row 14
$ws.Cells.Item(14, 3).Value = 3
$ws.Cells.Item(14, 6).Value = 121
$ws.Cells.Item(14, 8).Value = "living_rooms"
$ws.Cells.Item(14, 12).Value = "stimuli/img_bbs77.png"
$ws.Cells.Item(14, 13).Value = 31.64444444444445
$ws.Cells.Item(14, 14).Value = 21.26666666666667
$ws.Cells.Item(14, 15).Value = 26.45555555555556
$ws.Cells.Item(14, 16).Value = 45
$ws.Cells.Item(14, 17).Value = 2
$ws.Cells.Item(14, 18).Value = 2
$ws.Cells.Item(14, 19).Value = 2
# Row 15
$ws.Cells.Item(15, 3).Value = 3
$ws.Cells.Item(15, 6).Value = 122
$ws.Cells.Item(15, 12).Value = "stimuli/img_mdpr4.png"
$ws.Cells.Item(15, 13).Value = 74.04255319148936
$ws.Cells.Item(15, 14).Value = 54.70212765957447
$ws.Cells.Item(15, 15).Value = 64.37234042553192
$ws.Cells.Item(15, 16).Value = 47
$ws.Cells.Item(15, 17).Value = 6
$ws.Cells.Item(15, 18).Value = 6
$ws.Cells.Item(15, 19).Value = 6
# Row 16
$ws.Cells.Item(16, 3).Value = 3
$ws.Cells.Item(16, 6).Value = 123
$ws.Cells.Item(16, 8).Value = "living_rooms"
$ws.Cells.Item(16, 12).Value = "stimuli/img_6zz63.png"
$ws.Cells.Item(16, 13).Value = 87.66666666666667
$ws.Cells.Item(16, 14).Value = 70.6
$ws.Cells.Item(16, 15).Value = 79.13333333333333
$ws.Cells.Item(16, 16).Value = 45
$ws.Cells.Item(16, 17).Value = 9
$ws.Cells.Item(16, 18).Value = 10
$ws.Cells.Item(16, 19).Value = 10
# Row 17
$ws.Cells.Item(17, 3).Value = 3
$ws.Cells.Item(17, 6).Value = 124
$ws.Cells.Item(17, 8).Value = "living_rooms"
$ws.Cells.Item(17, 12).Value = "stimuli/img_eh0no.png"
$ws.Cells.Item(17, 13).Value = 53.66666666666666
$ws.Cells.Item(17, 14).Value = 36.02564102564103
$ws.Cells.Item(17, 15).Value = 44.84615384615385
$ws.Cells.Item(17, 16).Value = 39
$ws.Cells.Item(17, 17).Value = 3
$ws.Cells.Item(17, 18).Value = 3
$ws.Cells.Item(17, 19).Value = 3
# Row 18
$ws.Cells.Item(18, 3).Value = 3
$ws.Cells.Item(18, 6).Value = 125
$ws.Cells.Item(18, 8).Value = "living_rooms"
$ws.Cells.Item(18, 12).Value = "stimuli/img_abobq.png"
$ws.Cells.Item(18, 13).Value = 75.1842105263158
$ws.Cells.Item(18, 14).Value = 54.13157894736842
$ws.Cells.Item(18, 15).Value = 64.65789473684211
$ws.Cells.Item(18, 16).Value = 38
$ws.Cells.Item(18, 17).Value = 6
$ws.Cells.Item(18, 18).Value = 6
$ws.Cells.Item(18, 19).Value = 6
# Row 19
$ws.Cells.Item(19, 3).Value = 3
$ws.Cells.Item(19, 6).Value = 126
$ws.Cells.Item(19, 8).Value = "living_rooms"
$ws.Cells.Item(19, 12).Value = "stimuli/img_pey7u.png"
$ws.Cells.Item(19, 13).Value = 30.34883720930232
$ws.Cells.Item(19, 14).Value = 20.34883720930232
$ws.Cells.Item(19, 15).Value = 25.34883720930232
$ws.Cells.Item(19, 16).Value = 43
$ws.Cells.Item(19, 17).Value = 1
$ws.Cells.Item(19, 18).Value = 2
$ws.Cells.Item(19, 19).Value = 2
# Row 20
$ws.Cells.Item(20, 3).Value = 3
$ws.Cells.Item(20, 6).Value = 127
$ws.Cells.Item(20, 9).Value = "target"
$ws.Cells.Item(20, 11).Value = "j"
$ws.Cells.Item(20, 12).Value = "stimuli/img_wgkqa.png"
$ws.Cells.Item(20, 13).Value = 87.25581395348837
$ws.Cells.Item(20, 14).Value = 71.13953488372093
$ws.Cells.Item(20, 15).Value = 79.19767441860465
$ws.Cells.Item(20, 16).Value = 43
$ws.Cells.Item(20, 17).Value = 10
$ws.Cells.Item(20, 18).Value = 10
$ws.Cells.Item(20, 19).Value = 10
# Row 21
$ws.Cells.Item(21, 3).Value = 3
$ws.Cells.Item(21, 6).Value = 128
$ws.Cells.Item(21, 8).Value = "living_rooms"
$ws.Cells.Item(21, 12).Value = "stimuli/img_w8yhd.png"
$ws.Cells.Item(21, 13).Value = 55.74418604651163
$ws.Cells.Item(21, 14).Value = 38.90697674418605
$ws.Cells.Item(21, 15).Value = 47.32558139534883
$ws.Cells.Item(21, 16).Value = 43
$ws.Cells.Item(21, 17).Value = 4
$ws.Cells.Item(21, 18).Value = 4
$ws.Cells.Item(21, 19).Value = 4
# Row 22
$ws.Cells.Item(22, 3).Value = 3
$ws.Cells.Item(22, 6).Value = 129
$ws.Cells.Item(22, 9).Value = "target"
$ws.Cells.Item(22, 11).Value = "j"
$ws.Cells.Item(22, 12).Value = "stimuli/img_4o8l0.png"
$ws.Cells.Item(22, 13).Value = 46.02173913043478
$ws.Cells.Item(22, 14).Value = 31.45652173913043
$ws.Cells.Item(22, 15).Value = 38.73913043478261
$ws.Cells.Item(22, 16).Value = 46
$ws.Cells.Item(22, 17).Value = 3
$ws.Cells.Item(22, 18).Value = 3
$ws.Cells.Item(22, 19).Value = 3
# Row 23
$ws.Cells.Item(23, 3).Value = 3
$ws.Cells.Item(23, 6).Value = 130
$ws.Cells.Item(23, 8).Value = "living_rooms"
$ws.Cells.Item(23, 12).Value = "stimuli/img_bj99b.png"
$ws.Cells.Item(23, 13).Value = 82.79069767441861
$ws.Cells.Item(23, 14).Value = 65.46511627906976
$ws.Cells.Item(23, 15).Value = 74.12790697674419
$ws.Cells.Item(23, 16).Value = 43
$ws.Cells.Item(23, 17).Value = 8
$ws.Cells.Item(23, 18).Value = 8
$ws.Cells.Item(23, 19).Value = 8
# Row 24
$ws.Cells.Item(24, 3).Value = 3
$ws.Cells.Item(24, 6).Value = 131
$ws.Cells.Item(24, 8).Value = "living_rooms"
$ws.Cells.Item(24, 12).Value = "stimuli/img_xbtev.png"
$ws.Cells.Item(24, 13).Value = 13.68181818181818
$ws.Cells.Item(24, 14).Value = 8.568181818181818
$ws.Cells.Item(24, 15).Value = 11.125
$ws.Cells.Item(24, 16).Value = 44
$ws.Cells.Item(24, 17).Value = 1
$ws.Cells.Item(24, 18).Value = 1
$ws.Cells.Item(24, 19).Value = 1
# Row 25
$ws.Cells.Item(25, 3).Value = 3
$ws.Cells.Item(25, 6).Value = 132
$ws.Cells.Item(25, 8).Value = "living_rooms"
$ws.Cells.Item(25, 12).Value = "stimuli/img_xy930.png"
$ws.Cells.Item(25, 13).Value = 70.5952380952381
$ws.Cells.Item(25, 14).Value = 49.47619047619047
$ws.Cells.Item(25, 15).Value = 60.03571428571429
$ws.Cells.Item(25, 16).Value = 42
$ws.Cells.Item(25, 17).Value = 6
$ws.Cells.Item(25, 18).Value = 6
$ws.Cells.Item(25, 19).Value = 6
# Row 26
$ws.Cells.Item(26, 3).Value = 3
$ws.Cells.Item(26, 6).Value = 133
$ws.Cells.Item(26, 8).Value = "bedrooms"
$ws.Cells.Item(26, 9).Value = "distractor"
$ws.Cells.Item(26, 11).Value = "f"
$ws.Cells.Item(26, 12).Value = "stimuli/img_5p2ql.png"
$ws.Cells.Item(26, 13).Value = 89.19565217391305
$ws.Cells.Item(26, 14).Value = 72.52173913043478
$ws.Cells.Item(26, 15).Value = 80.8586956521739
$ws.Cells.Item(26, 16).Value = 46
$ws.Cells.Item(26, 17).Value = 10
$ws.Cells.Item(26, 18).Value = 10
$ws.Cells.Item(26, 19).Value = 10
# Row 27
$ws.Cells.Item(27, 3).Value = 3
$ws.Cells.Item(27, 6).Value = 134
$ws.Cells.Item(27, 8).Value = "living_rooms"
$ws.Cells.Item(27, 12).Value = "stimuli/img_di6f0.png"
$ws.Cells.Item(27, 13).Value = 94.04347826086956
$ws.Cells.Item(27, 14).Value = 83.34782608695652
$ws.Cells.Item(27, 15).Value = 88.69565217391303
$ws.Cells.Item(27, 16).Value = 46

# Remove the last row (row 28), which no longer exists in the updated dataset
$ws.Rows.Item(28).Delete()

Write-Host "Done. Changed cells:" 290
